$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each row update for the crypto price/volume refresh.
# D column values are forced to Text format so that strings such as
# "1.00" or "58.326.37" are preserved exactly instead of being
# reinterpreted by Excel as numbers (which would drop trailing zeros
# or introduce floating point artifacts).

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.326.37"
$ws.Range("E2").Value = "  +2.14%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.149.24"
$ws.Range("E3").Value = "  +2.62%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.87"
$ws.Range("E5").Value = "  +3.01%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.80"
$ws.Range("E6").Value = "  +3.23%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"

# Row 8: XRP
$ws.Range("E8").Value = "  +8.85%  "

# Row 9: Toncoin
$ws.Range("E9").Value = "  +0.39%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +3.23%  "

# Row 11: Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.422"
$ws.Range("E11").Value = "  +5.05%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +2.92%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.693.40"
$ws.Range("E13").Value = "  +2.70%  "

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.07"
$ws.Range("E14").Value = "  +3.57%  "

# Row 15: ShibaInu
$ws.Range("E15").Value = "  +5.76%  "

# Row 16: WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.371.54"
$ws.Range("E16").Value = "  +2.12%  "

# Row 17: WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.154.26"
$ws.Range("E17").Value = "  +3.01%  "

# Row 18: Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").Value = "  +6.25%  "

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +4.62%  "

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("E20").Value = "  +4.83%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.46"
$ws.Range("E21").Value = "  +8.27%  "

# Row 22: Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23: LEO
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.76"
$ws.Range("E23").Value = "  -0.42%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.57"
$ws.Range("E24").Value = "  +2.15%  "

# Row 25: Polygon
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.517"
$ws.Range("E25").Value = "  +3.90%  "

# Row 26: Kaspa
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +1.71%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.28%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +12.45%  "

# Row 29: PEPE
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0881"
$ws.Range("E29").Value = "  +2.69%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  +2.66%  "

# Row 31: RenderToken
$ws.Range("E31").Value = "  +6.85%  "

# Row 32: EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.81"
$ws.Range("E32").Value = "  +4.38%  "

# Row 33: NEARProtocol
$ws.Range("E33").Value = "  +7.55%  "

# Row 34: Fetch.AI
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.17"
$ws.Range("E34").Value = "  +4.78%  "

# Row 35: Monero
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.29"
$ws.Range("E35").Value = "  +1.55%  "

# Row 36: Aptos
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.25"
$ws.Range("E36").Value = "  +4.69%  "

# Row 37: ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.35"
$ws.Range("E37").Value = "  +9.94%  "

# Row 38: EnergySwap
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.49"
$ws.Range("E38").Value = "  +0.92%  "

# Row 39: Stacks
$ws.Range("E39").Value = "  +7.85%  "

# Row 40: Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.649.60"
$ws.Range("E40").Value = "  +10.10%  "

# Row 41: Hedera
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0680"
$ws.Range("E41").Value = "  +3.97%  "

# Row 42: Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.23"
$ws.Range("E42").Value = "  +5.14%  "

# Row 43: OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.63"
$ws.Range("E43").Value = "  +6.02%  "

# Row 44: Mantle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.703"
$ws.Range("E44").Value = "  +1.74%  "

# Row 45: VeChain
$ws.Range("E45").Value = "  +5.97%  "

# Row 47: Stellar
$ws.Range("E47").Value = "  +11.71%  "

# Row 48: Cosmos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.23"
$ws.Range("E48").Value = "  +4.61%  "

# Row 49: ONDO
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.978"
$ws.Range("E49").Value = "  +5.13%  "

# Row 50: InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.27"
$ws.Range("E50").Value = "  +4.24%  "

# Row 51: SuiNetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.753"
$ws.Range("E51").Value = "  +0.57%  "
